$wb = $excel.ActiveWorkbook

# Fill in column C ("Associated Variable with Margin of Error") for every
# data row on the "Blank Template" sheet with the matching "<variable>_moe"
# name derived from column B ("Element or value display name").
$ws1 = $wb.Worksheets.Item("Blank Template")
for ($r = 2; $r -le 43; $r++) {
    $bVal = $ws1.Cells.Item($r, 2).Value2
    $ws1.Cells.Item($r, 3).Value = $bVal + "_moe"
}

# Update the active selection to match the new state of the workbook.
[void]$ws1.Range("D4").Select()

# Remove the separate "Description of Fields" instructions sheet - its
# content is no longer needed.
$excel.DisplayAlerts = $false
$ws2 = $wb.Worksheets.Item("Description of Fields")
[void]$ws2.Delete()
$excel.DisplayAlerts = $true

# Rename the remaining sheet to reflect its new purpose.
$ws1.Name = "List of Variables"
